$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 356, pushing the existing rows 356-397
# down to 357-398 (this also grows the sheet dimension to A1:R398).
$ws.Rows.Item(356).Insert()

# Populate the newly inserted row 356 with a new weekly price record for
# "Feria Lagunitas de Puerto Montt" / Ajo, reusing the fixed attributes
# that are constant across this subset and filling in the new date /
# volume / price data.
$ws.Cells.Item(356, 1).Value  = 4
$ws.Cells.Item(356, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(356, 3).Value  = "Los Lagos"
$ws.Cells.Item(356, 4).Value  = 44946
$ws.Cells.Item(356, 5).Value  = 10
$ws.Cells.Item(356, 6).Value  = 100112003
$ws.Cells.Item(356, 7).Value  = "Ajo"
$ws.Cells.Item(356, 8).Value  = "Chino"
$ws.Cells.Item(356, 9).Value  = "Primera"
$ws.Cells.Item(356, 10).Value = 240
$ws.Cells.Item(356, 11).Value = 19000
$ws.Cells.Item(356, 12).Value = 20000
$ws.Cells.Item(356, 13).Value = 19500
$ws.Cells.Item(356, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(356, 15).Value = "China"
$ws.Cells.Item(356, 16).Value = 1950
$ws.Cells.Item(356, 17).Value = 10
$ws.Cells.Item(356, 18).Value = "Hortaliza"
